$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule "R30" (row 10) lower bound ("From", column C) is updated from 18 to 1.
$ws.Range("C10").Value = 1
